$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2: "Passed" -> "Failed:Button Disabled For This Manager"
$ws.Range("D2").Value = "Failed:Button Disabled For This Manager"

# Remove D3 content entirely (cell no longer present in sheet data)
$ws.Range("D3").ClearContents()

# Move the active selection to D2 (previously D5)
$ws.Range("D2").Select()
